# Applies the scheduled-runner update to the Leve profit tables (H:N)
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 110898.5
$ws.Range("I6").Value = 75.5
$ws.Range("K6").Value = 226.5
$ws.Range("M6").Value = -114.5
# Row 8
$ws.Range("H8").Value = 1208.6666
$ws.Range("I8").Value = 38
$ws.Range("J8").Value = 3550
$ws.Range("K8").Value = 114
$ws.Range("L8").Value = 10650
$ws.Range("M8").Value = 25
$ws.Range("N8").Value = -10928
# Row 21
$ws.Range("H21").Value = 8500
$ws.Range("I21").Value = 7000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -6532
$ws.Range("N21").Value = -10936
# Row 23
$ws.Range("H23").Value = 8500
$ws.Range("I23").Value = 7000
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -6766
$ws.Range("N23").Value = -10468
# Row 70
$ws.Range("H70").Value = 4625
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 6083.3335
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 18250.0005
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -18790.0005
# Row 73
$ws.Range("H73").Value = 4625
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 6083.3335
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 18250.0005
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -20122.0005
# Row 111
$ws.Range("H111").Value = 1700
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1700
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 5100
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -11234
# Row 113
$ws.Range("H113").Value = 4030.5
$ws.Range("I113").Value = 3975.625
$ws.Range("K113").Value = 3975.625
$ws.Range("M113").Value = -721.625
# Row 116
$ws.Range("H116").Value = 2028
$ws.Range("I116").Value = 1613.3334
$ws.Range("J116").Value = 2650
$ws.Range("K116").Value = 1613.3334
$ws.Range("L116").Value = 2650
$ws.Range("M116").Value = 1828.6666
$ws.Range("N116").Value = -9534
# Row 127
$ws.Range("H127").Value = 1051.5555
$ws.Range("I127").Value = 517.63635
$ws.Range("J127").Value = 1890.5714
$ws.Range("K127").Value = 1552.90905
$ws.Range("L127").Value = 5671.7142
$ws.Range("M127").Value = 3407.09095
$ws.Range("N127").Value = -15591.7142
# Row 129
$ws.Range("H129").Value = 991.7353000000001
$ws.Range("I129").Value = 627.5
$ws.Range("J129").Value = 1014.5
$ws.Range("K129").Value = 1882.5
$ws.Range("L129").Value = 3043.5
$ws.Range("M129").Value = 3117.5
$ws.Range("N129").Value = -13043.5
# Row 132
$ws.Range("H132").Value = 1022624.75
$ws.Range("I132").Value = 1617.9302
$ws.Range("J132").Value = 9803283
$ws.Range("K132").Value = 4853.7906
$ws.Range("L132").Value = 29409849
$ws.Range("M132").Value = -2323.7906
$ws.Range("N132").Value = -29414909
# Row 135
$ws.Range("H135").Value = 25461.854
$ws.Range("I135").Value = 30122.766
$ws.Range("J135").Value = 2823.1428
$ws.Range("K135").Value = 271104.894
$ws.Range("L135").Value = 25408.2852
$ws.Range("M135").Value = -268569.894
$ws.Range("N135").Value = -30478.2852
# Row 138
$ws.Range("H138").Value = 2417937
$ws.Range("I138").Value = 805.7241
$ws.Range("J138").Value = 4170357.2
$ws.Range("K138").Value = 2417.1723
$ws.Range("L138").Value = 12511071.6
$ws.Range("M138").Value = 2722.8277
$ws.Range("N138").Value = -12521351.6

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 45546428
$ws.Range("I61").Value = 52685184
$ws.Range("J61").Value = 334304.66
$ws.Range("K61").Value = 52685184
$ws.Range("L61").Value = 334304.66
$ws.Range("M61").Value = -52684972
$ws.Range("N61").Value = -334728.66
# Row 63
$ws.Range("H63").Value = 2984.1667
$ws.Range("I63").Value = 2984.1667
$ws.Range("K63").Value = 2984.1667
$ws.Range("M63").Value = -2298.1667
# Row 66
$ws.Range("H66").Value = 2984.1667
$ws.Range("I66").Value = 2984.1667
$ws.Range("K66").Value = 14920.8335
$ws.Range("M66").Value = -11488.8335
# Row 74
$ws.Range("H74").Value = 4422691
$ws.Range("I74").Value = 5578867
$ws.Range("J74").Value = 87031.664
$ws.Range("K74").Value = 5578867
$ws.Range("L74").Value = 87031.664
$ws.Range("M74").Value = -5577993
$ws.Range("N74").Value = -88779.664
# Row 77
$ws.Range("H77").Value = 4422691
$ws.Range("I77").Value = 5578867
$ws.Range("J77").Value = 87031.664
$ws.Range("K77").Value = 27894335
$ws.Range("L77").Value = 435158.32
$ws.Range("M77").Value = -27889967
$ws.Range("N77").Value = -443894.32
# Row 132
$ws.Range("H132").Value = 111947
$ws.Range("I132").Value = 84499.586
$ws.Range("K132").Value = 253498.758
$ws.Range("M132").Value = -250968.758
# Row 136
$ws.Range("H136").Value = 45546428
$ws.Range("I136").Value = 52685184
$ws.Range("J136").Value = 334304.66
$ws.Range("K136").Value = 158055552
$ws.Range("L136").Value = 1002913.98
$ws.Range("M136").Value = -158053002
$ws.Range("N136").Value = -1008013.98

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1078.5143
$ws.Range("I99").Value = 1086.95
$ws.Range("K99").Value = 1086.95
$ws.Range("M99").Value = 411.05
# Row 134
$ws.Range("H134").Value = 3223.2856
$ws.Range("I134").Value = 2828.2046
$ws.Range("J134").Value = 6700
$ws.Range("K134").Value = 8484.613799999999
$ws.Range("L134").Value = 20100
$ws.Range("M134").Value = -5949.613799999999
$ws.Range("N134").Value = -25170

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 2857.4375
$ws.Range("I99").Value = 2187.12
$ws.Range("J99").Value = 5251.4287
$ws.Range("K99").Value = 2187.12
$ws.Range("L99").Value = 5251.4287
$ws.Range("M99").Value = -689.1199999999999
$ws.Range("N99").Value = -8247.4287
# Row 126
$ws.Range("H126").Value = 2857.4375
$ws.Range("I126").Value = 2187.12
$ws.Range("J126").Value = 5251.4287
$ws.Range("K126").Value = 6561.36
$ws.Range("L126").Value = 15754.2861
$ws.Range("M126").Value = -4091.36
$ws.Range("N126").Value = -20694.2861

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 2227.7646
$ws.Range("I10").Value = 1089.6666
$ws.Range("J10").Value = 4959.2
$ws.Range("K10").Value = 3268.9998
$ws.Range("L10").Value = 14877.6
$ws.Range("M10").Value = -3129.9998
$ws.Range("N10").Value = -15155.6
# Row 16
$ws.Range("H16").Value = 1160
$ws.Range("J16").Value = 1375
$ws.Range("L16").Value = 4125
$ws.Range("N16").Value = -4471
# Row 80
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 12000
$ws.Range("N80").Value = -13872
# Row 83
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 36000
$ws.Range("N83").Value = -45360
# Row 92
$ws.Range("H92").Value = 954.1053000000001
$ws.Range("I92").Value = 949.0769
$ws.Range("J92").Value = 965
$ws.Range("K92").Value = 2847.2307
$ws.Range("L92").Value = 2895
$ws.Range("M92").Value = -1599.2307
$ws.Range("N92").Value = -5391
# Row 131
$ws.Range("H131").Value = 1304
$ws.Range("J131").Value = 1402.0513
$ws.Range("L131").Value = 4206.1539
$ws.Range("N131").Value = -14286.1539

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1925
$ws.Range("I113").Value = 1750
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1750
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 420
$ws.Range("N113").Value = -6440
# Row 132
$ws.Range("H132").Value = 113069.22
$ws.Range("I132").Value = 68349.87
$ws.Range("K132").Value = 205049.61
$ws.Range("M132").Value = -202519.61
# Row 141
$ws.Range("H141").Value = 29397.143
$ws.Range("J141").Value = 29397.143
$ws.Range("L141").Value = 29397.143
$ws.Range("N141").Value = -39757.143

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2837.25
$ws.Range("J7").Value = 2856.8572
$ws.Range("L7").Value = 2856.8572
$ws.Range("N7").Value = -3080.8572
# Row 93
$ws.Range("H93").Value = 2375
$ws.Range("J93").Value = 2250
$ws.Range("L93").Value = 2250
$ws.Range("N93").Value = -4746
# Row 122
$ws.Range("H122").Value = 3802.6667
$ws.Range("I122").Value = 3172.4
$ws.Range("J122").Value = 4252.857
$ws.Range("K122").Value = 9517.200000000001
$ws.Range("L122").Value = 12758.571
$ws.Range("M122").Value = -7067.200000000001
$ws.Range("N122").Value = -17658.571
# Row 126
$ws.Range("H126").Value = 2837.25
$ws.Range("J126").Value = 2856.8572
$ws.Range("L126").Value = 8570.571599999999
$ws.Range("N126").Value = -13510.5716
# Row 132
$ws.Range("H132").Value = 65962.03999999999
$ws.Range("I132").Value = 32440.354
$ws.Range("J132").Value = 147371.86
$ws.Range("K132").Value = 97321.06200000001
$ws.Range("L132").Value = 442115.58
$ws.Range("M132").Value = -94791.06200000001
$ws.Range("N132").Value = -447175.58
# Row 136
$ws.Range("H136").Value = 58227.723
$ws.Range("I136").Value = 32568.688
$ws.Range("K136").Value = 97706.064
$ws.Range("M136").Value = -95156.064
